# Update the Haba (Vega Modelo de Temuco) weekly price records.
# The data rows (3-22) are re-shuffled: each row now shows the values
# that, in the previous snapshot, belonged to a different day in the week.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 4).Value = 44356
$ws.Cells.Item(3, 10).Value = 30
$ws.Cells.Item(3, 11).Value = 14000
$ws.Cells.Item(3, 12).Value = 14000
$ws.Cells.Item(3, 13).Value = 14000
$ws.Cells.Item(3, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 16).Value = 560

$ws.Cells.Item(4, 4).Value = 44160
$ws.Cells.Item(4, 10).Value = 30
$ws.Cells.Item(4, 11).Value = 8000
$ws.Cells.Item(4, 12).Value = 8000
$ws.Cells.Item(4, 13).Value = 8000
$ws.Cells.Item(4, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(4, 16).Value = 320

$ws.Cells.Item(5, 4).Value = 44354
$ws.Cells.Item(5, 10).Value = 80
$ws.Cells.Item(5, 11).Value = 16000
$ws.Cells.Item(5, 12).Value = 16000
$ws.Cells.Item(5, 13).Value = 16000
$ws.Cells.Item(5, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 16).Value = 640

$ws.Cells.Item(6, 4).Value = 44162
$ws.Cells.Item(6, 10).Value = 260
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 12).Value = 8000
$ws.Cells.Item(6, 13).Value = 7462
$ws.Cells.Item(6, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(6, 16).Value = 298

$ws.Cells.Item(7, 4).Value = 44181
$ws.Cells.Item(7, 10).Value = 55
$ws.Cells.Item(7, 11).Value = 14000
$ws.Cells.Item(7, 12).Value = 14000
$ws.Cells.Item(7, 13).Value = 14000
$ws.Cells.Item(7, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(7, 16).Value = 560

$ws.Cells.Item(8, 4).Value = 44355
$ws.Cells.Item(8, 10).Value = 20
$ws.Cells.Item(8, 11).Value = 16000
$ws.Cells.Item(8, 12).Value = 16000
$ws.Cells.Item(8, 13).Value = 16000
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 640

$ws.Cells.Item(9, 4).Value = 44159
$ws.Cells.Item(9, 10).Value = 50
$ws.Cells.Item(9, 11).Value = 8000
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 8000
$ws.Cells.Item(9, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(9, 16).Value = 320

$ws.Cells.Item(10, 4).Value = 44159
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(10, 11).Value = 8000
$ws.Cells.Item(10, 12).Value = 8000
$ws.Cells.Item(10, 13).Value = 8000
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 320

$ws.Cells.Item(11, 4).Value = 44186
$ws.Cells.Item(11, 10).Value = 30
$ws.Cells.Item(11, 11).Value = 14000
$ws.Cells.Item(11, 12).Value = 14000
$ws.Cells.Item(11, 13).Value = 14000
$ws.Cells.Item(11, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(11, 16).Value = 560

$ws.Cells.Item(12, 4).Value = 44371
$ws.Cells.Item(12, 10).Value = 40
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 15000
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 16).Value = 600

$ws.Cells.Item(13, 4).Value = 44392
$ws.Cells.Item(13, 10).Value = 55
$ws.Cells.Item(13, 11).Value = 17000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = 17455
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 698

$ws.Cells.Item(14, 4).Value = 44168
$ws.Cells.Item(14, 10).Value = 120
$ws.Cells.Item(14, 11).Value = 7000
$ws.Cells.Item(14, 12).Value = 8000
$ws.Cells.Item(14, 13).Value = 7458
$ws.Cells.Item(14, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(14, 16).Value = 298

$ws.Cells.Item(15, 4).Value = 44434
$ws.Cells.Item(15, 10).Value = 50
$ws.Cells.Item(15, 11).Value = 15000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 15000
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 600

$ws.Cells.Item(16, 4).Value = 44210
$ws.Cells.Item(16, 10).Value = 110
$ws.Cells.Item(16, 11).Value = 16000
$ws.Cells.Item(16, 12).Value = 16000
$ws.Cells.Item(16, 13).Value = 16000
$ws.Cells.Item(16, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(16, 16).Value = 640

$ws.Cells.Item(17, 4).Value = 44427
$ws.Cells.Item(17, 10).Value = 30
$ws.Cells.Item(17, 11).Value = 15000
$ws.Cells.Item(17, 12).Value = 15000
$ws.Cells.Item(17, 13).Value = 15000
$ws.Cells.Item(17, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 16).Value = 600

$ws.Cells.Item(18, 4).Value = 44176
$ws.Cells.Item(18, 10).Value = 20
$ws.Cells.Item(18, 11).Value = 11000
$ws.Cells.Item(18, 12).Value = 11000
$ws.Cells.Item(18, 13).Value = 11000
$ws.Cells.Item(18, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(18, 16).Value = 440

$ws.Cells.Item(19, 4).Value = 44161
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 7000
$ws.Cells.Item(19, 12).Value = 7000
$ws.Cells.Item(19, 13).Value = 7000
$ws.Cells.Item(19, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(19, 16).Value = 280

$ws.Cells.Item(20, 4).Value = 44435
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 15000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 13).Value = 15000
$ws.Cells.Item(20, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 16).Value = 600

$ws.Cells.Item(21, 4).Value = 44175
$ws.Cells.Item(21, 10).Value = 50
$ws.Cells.Item(21, 11).Value = 8000
$ws.Cells.Item(21, 12).Value = 8000
$ws.Cells.Item(21, 13).Value = 8000
$ws.Cells.Item(21, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(21, 16).Value = 320

$ws.Cells.Item(22, 4).Value = 44167
$ws.Cells.Item(22, 10).Value = 95
$ws.Cells.Item(22, 11).Value = 7000
$ws.Cells.Item(22, 12).Value = 7000
$ws.Cells.Item(22, 13).Value = 7000
$ws.Cells.Item(22, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(22, 16).Value = 280

